$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.739.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.448.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.28%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.449.49"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.35%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.473"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.49%  "
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.031.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.122"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.68%  "
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.453.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.821.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.51%  "
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "385.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.564"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.590.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.08%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.27%  "
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "72.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.177"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -13.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "24.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.28%  "
$ws.Range("E39").Value = "  +1.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "166.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0789"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.69%  "
$ws.Range("E43").Value = "  +2.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.68%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "42.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.52%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.616.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.36%  "
